# Update the "想去人数" (wish-to-attend) counts on the 展览 and 全部类型
# sheets. Both sheets mirror the same listing, so the same five rows are
# updated on each.

$wb = $excel.ActiveWorkbook

$updates = @{
    6  = 4587
    9  = 1341
    12 = 970
    14 = 541
    15 = 57
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
